$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 33, shifting the existing rows
# 33-36 down to 34-37 (their data/formatting move with them automatically).
$ws.Rows.Item(33).Insert()

# Populate the newly inserted row 33 with the new weekly record.
$ws.Range("A33").Value = 5
$ws.Range("B33").Value = "Macroferia Regional de Talca"
$ws.Range("C33").Value = "Maule"
$ws.Range("D33").Value = 44746
$ws.Range("E33").Value = 7
$ws.Range("F33").Value = 100112043
$ws.Range("G33").Value = "Pepino dulce"
$ws.Range("H33").Value = "Cultivar IV Región"
$ws.Range("I33").Value = "Primera"
$ws.Range("J33").Value = 500
$ws.Range("K33").Value = 15000
$ws.Range("L33").Value = 15000
$ws.Range("M33").Value = 15000
$ws.Range("N33").Value = "`$/bandeja 18 kilos"
$ws.Range("O33").Value = "Provincia de Limarí"
$ws.Range("P33").Value = 833
$ws.Range("Q33").Value = 18
$ws.Range("R33").Value = "Hortaliza"
